# Apply odds updates to Sheet1, as described in the commit diff
# ("Atualizando o arquivo XLSX").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Platense vs Rosario Central)
$ws.Range("G2").Value = 2.4
$ws.Range("H2").Value = 2.7
$ws.Range("I2").Value = 3.8
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 1.9
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = 1.36
$ws.Range("X2").Value = 34
$ws.Range("Z2").Value = 4
$ws.Range("AI2").Value = 41

# Row 3 (Ceara vs Vasco)
$ws.Range("G3").Value = 2.25
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 3.4
$ws.Range("N3").Value = 2.35
$ws.Range("O3").Value = 1.57
$ws.Range("Z3").Value = 7
$ws.Range("AF3").Value = 15

# Row 9 (Pereira vs Dep. Cali)
$ws.Range("G9").Value = 2.45
$ws.Range("H9").Value = 2.9
$ws.Range("I9").Value = 3.2
$ws.Range("L9").Value = 1.53
$ws.Range("M9").Value = 2.38
$ws.Range("N9").Value = 2.7
$ws.Range("O9").Value = 1.44
$ws.Range("U9").Value = 10
$ws.Range("W9").Value = 23
$ws.Range("X9").Value = 23
$ws.Range("Z9").Value = 6.5
$ws.Range("AE9").Value = 7.5
$ws.Range("AF9").Value = 15
$ws.Range("AG9").Value = 13
$ws.Range("AH9").Value = 34

# Row 10 (Deportes Tolima vs Junior)
$ws.Range("I10").Value = 4.33
$ws.Range("J10").Value = 1.1
$ws.Range("K10").Value = 7
$ws.Range("AG10").Value = 15

# Row 12 (Herediano vs AD Santos)
$ws.Range("G12").Value = 1.42
$ws.Range("L12").Value = 1.22
$ws.Range("M12").Value = 4

# Row 13 (Puntarenas FC vs Saprissa)
$ws.Range("I13").Value = 1.95

# Row 14 (Aucas vs U. Catolica)
$ws.Range("G14").Value = 2.25
$ws.Range("H14").Value = 3.3
$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 1.05
$ws.Range("K14").Value = 11
$ws.Range("N14").Value = 1.95
$ws.Range("O14").Value = 1.85
$ws.Range("T14").Value = 8
$ws.Range("U14").Value = 11
$ws.Range("V14").Value = 9.5
$ws.Range("W14").Value = 21
$ws.Range("X14").Value = 19
$ws.Range("Y14").Value = 26
$ws.Range("AE14").Value = 10
$ws.Range("AF14").Value = 15
$ws.Range("AG14").Value = 11
$ws.Range("AH14").Value = 34
$ws.Range("AI14").Value = 23
$ws.Range("AJ14").Value = 34

# Row 18 (Pachuca vs Tigres UANL)
$ws.Range("G18").Value = 1.91
$ws.Range("L18").Value = 1.3
$ws.Range("M18").Value = 3.4
$ws.Range("O18").Value = 1.75
$ws.Range("P18").Value = 1.41
$ws.Range("Q18").Value = 2.62

# Row 19 (Guadalajara Chivas vs Puebla)
$ws.Range("I19").Value = 7
$ws.Range("N19").Value = 1.9
$ws.Range("O19").Value = 1.95
$ws.Range("T19").Value = 6
